# "Add files via upload" — the author re-uploaded a trimmed-down version of
# the workbook: only the "Datos Generales" sheet survives; the "Evaluación"
# and "Resumen" sheets (which just duplicated/summarised the same scoring
# table) are removed entirely.

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete this sheet" prompt blocking
# the automation, mirroring typical Excel COM scripting practice.
$excel.DisplayAlerts = $false

$wb.Worksheets("Evaluación").Delete() | Out-Null
$wb.Worksheets("Resumen").Delete() | Out-Null

$excel.DisplayAlerts = $true
